# Rename the inline logo pictures in the document's headers/footers.
#
#   headers (BTec_Logo-Orange): image1.jpg -> image2.jpg
#   footers (PearsonLogo.png):  image2.png -> image1.png
#
# InlineShape has no writable "Name" in the Word object model, so each
# picture is temporarily converted to a floating Shape (which does expose
# Name), renamed, then converted back to an inline picture in place.

$d = $word.ActiveDocument
$sec = $d.Sections(1)

for ($i = 1; $i -le $sec.Headers.Count; $i++) {
    $hf = $sec.Headers($i)
    if (-not $hf.Exists) { continue }

    $shapeCount = $hf.Range.InlineShapes.Count
    for ($j = 1; $j -le $shapeCount; $j++) {
        $inline = $hf.Range.InlineShapes($j)
        $shape = $inline.ConvertToShape()

        $newName = $null
        if ($shape.Name -eq "image1.jpg") { $newName = "image2.jpg" }
        elseif ($shape.Name -eq "image2.png") { $newName = "image1.png" }

        if ($newName -ne $null) {
            $shape.Name = $newName
        }

        $shape.ConvertToInlineShape() | Out-Null
    }
}

for ($i = 1; $i -le $sec.Footers.Count; $i++) {
    $hf = $sec.Footers($i)
    if (-not $hf.Exists) { continue }

    $shapeCount = $hf.Range.InlineShapes.Count
    for ($j = 1; $j -le $shapeCount; $j++) {
        $inline = $hf.Range.InlineShapes($j)
        $shape = $inline.ConvertToShape()

        $newName = $null
        if ($shape.Name -eq "image1.jpg") { $newName = "image2.jpg" }
        elseif ($shape.Name -eq "image2.png") { $newName = "image1.png" }

        if ($newName -ne $null) {
            $shape.Name = $newName
        }

        $shape.ConvertToInlineShape() | Out-Null
    }
}

Write-Output "Renamed header/footer logo pictures."
